$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates
$ws.Cells.Item(2, 4).Value = "34.953.55"
$ws.Cells.Item(3, 4).Value = "1.845.69"
$ws.Cells.Item(5, 4).Value = "'232.56"
$ws.Cells.Item(8, 4).Value = "'40.80"
$ws.Cells.Item(9, 4).Value = "'0.329"
$ws.Cells.Item(11, 4).Value = "'0.0985"
$ws.Cells.Item(12, 4).Value = "2.113.95"
$ws.Cells.Item(13, 4).Value = "'11.44"
$ws.Cells.Item(14, 4).Value = "1.850.83"
$ws.Cells.Item(16, 4).Value = "'4.65"
$ws.Cells.Item(17, 4).Value = "34.997.07"
$ws.Cells.Item(18, 4).Value = "'69.97"
$ws.Cells.Item(20, 4).Value = "'240.34"
$ws.Cells.Item(24, 4).Value = "'2.25"
$ws.Cells.Item(25, 4).Value = "'172.79"
$ws.Cells.Item(26, 4).Value = "'7.84"
$ws.Cells.Item(27, 4).Value = "'17.46"
$ws.Cells.Item(28, 4).Value = "'0.123"
$ws.Cells.Item(31, 4).Value = "'0.0553"
$ws.Cells.Item(35, 4).Value = "'1.95"
$ws.Cells.Item(36, 4).Value = "'0.749"
$ws.Cells.Item(39, 4).Value = "'90.02"
$ws.Cells.Item(40, 4).Value = "1.347.49"
$ws.Cells.Item(42, 4).Value = "'14.67"
$ws.Cells.Item(47, 4).Value = "'6.31"
$ws.Cells.Item(48, 4).Value = "2.030.64"
$ws.Cells.Item(51, 4).Value = "'0.0670"

# Column E (Volume 1h) updates
$ws.Cells.Item(2, 5).Value = "  +0.21%  "
$ws.Cells.Item(3, 5).Value = "  +2.03%  "
$ws.Cells.Item(4, 5).Value = "  +0.12%  "
$ws.Cells.Item(5, 5).Value = "  +0.43%  "
$ws.Cells.Item(6, 5).Value = "  +2.16%  "
$ws.Cells.Item(7, 5).Value = "  +0.18%  "
$ws.Cells.Item(8, 5).Value = "  +4.50%  "
$ws.Cells.Item(9, 5).Value = "  +3.65%  "
$ws.Cells.Item(10, 5).Value = "  +1.83%  "
$ws.Cells.Item(11, 5).Value = "  -0.92%  "
$ws.Cells.Item(12, 5).Value = "  +2.12%  "
$ws.Cells.Item(13, 5).Value = "  +4.97%  "
$ws.Cells.Item(14, 5).Value = "  +3.81%  "
$ws.Cells.Item(15, 5).Value = "  +2.13%  "
$ws.Cells.Item(16, 5).Value = "  +2.18%  "
$ws.Cells.Item(18, 5).Value = "  +1.44%  "
$ws.Cells.Item(19, 5).Value = "  +1.40%  "
$ws.Cells.Item(20, 5).Value = "  +1.10%  "
$ws.Cells.Item(21, 5).Value = "  +3.83%  "
$ws.Cells.Item(22, 5).Value = "  +2.53%  "
$ws.Cells.Item(23, 5).Value = "  +0.20%  "
$ws.Cells.Item(24, 5).Value = "  +0.88%  "
$ws.Cells.Item(25, 5).Value = "  +0.55%  "
$ws.Cells.Item(26, 5).Value = "  +0.95%  "
$ws.Cells.Item(27, 5).Value = "  +1.66%  "
$ws.Cells.Item(28, 5).Value = "  +3.61%  "
$ws.Cells.Item(29, 5).Value = "  +5.40%  "
$ws.Cells.Item(30, 5).Value = "  +0.13%  "
$ws.Cells.Item(31, 5).Value = "  +1.17%  "
$ws.Cells.Item(32, 5).Value = "  +0.09%  "
$ws.Cells.Item(33, 5).Value = "  +0.97%  "
$ws.Cells.Item(34, 5).Value = "  +23.38%  "
$ws.Cells.Item(35, 5).Value = "  +11.13%  "
$ws.Cells.Item(36, 5).Value = "  +10.16%  "
$ws.Cells.Item(37, 5).Value = "  +6.13%  "
$ws.Cells.Item(38, 5).Value = "  +12.05%  "
$ws.Cells.Item(39, 5).Value = "  -1.18%  "
$ws.Cells.Item(40, 5).Value = "  +3.28%  "
$ws.Cells.Item(41, 5).Value = "  +2.77%  "
$ws.Cells.Item(42, 5).Value = "  +2.48%  "
$ws.Cells.Item(43, 5).Value = "  +3.29%  "
$ws.Cells.Item(44, 5).Value = "  -1.60%  "
$ws.Cells.Item(45, 5).Value = "  +2.39%  "
$ws.Cells.Item(46, 5).Value = "  +4.22%  "
$ws.Cells.Item(47, 5).Value = "  +3.08%  "
$ws.Cells.Item(48, 5).Value = "  +2.03%  "
$ws.Cells.Item(49, 5).Value = "  +18.29%  "
$ws.Cells.Item(50, 5).Value = "  +0.25%  "
$ws.Cells.Item(51, 5).Value = "  -0.56%  "
